# Update computed price/profit figures across the Leve profitability sheets.
# (ALC, ARM, CRP, CUL, GSM, LTW, WVR - BSM is unchanged)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6527.625
$ws.Range("I28").Value = 7969.077
$ws.Range("K28").Value = 7969.077
$ws.Range("M28").Value = -7484.077

$ws.Range("H70").Value = 125003670
$ws.Range("J70").Value = 166671150
$ws.Range("L70").Value = 500013450
$ws.Range("N70").Value = -500013990

$ws.Range("H73").Value = 125003670
$ws.Range("J73").Value = 166671150
$ws.Range("L73").Value = 500013450
$ws.Range("N73").Value = -500015322

$ws.Range("H111").Value = 2350.5454
$ws.Range("J111").Value = 2596
$ws.Range("L111").Value = 7788
$ws.Range("N111").Value = -13922

$ws.Range("H136").Value = 107960
$ws.Range("J136").Value = 107960
$ws.Range("L136").Value = 107960
$ws.Range("N136").Value = -118160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 350000
$ws.Range("J34").Value = 500000
$ws.Range("L34").Value = 500000
$ws.Range("N34").Value = -500542

$ws.Range("H80").Value = 43666.668
$ws.Range("J80").Value = 48400
$ws.Range("L80").Value = 48400
$ws.Range("N80").Value = -50396

$ws.Range("H83").Value = 43666.668
$ws.Range("J83").Value = 48400
$ws.Range("L83").Value = 145200
$ws.Range("N83").Value = -155184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 25083
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H58").Value = 1324
$ws.Range("I58").Value = 1286.5
$ws.Range("J58").Value = 1499
$ws.Range("K58").Value = 1286.5
$ws.Range("L58").Value = 1499
$ws.Range("M58").Value = -1083.5
$ws.Range("N58").Value = -1905

$ws.Range("H132").Value = 2131.2778
$ws.Range("I132").Value = 2138.1765
$ws.Range("J132").Value = 2014
$ws.Range("K132").Value = 6414.529500000001
$ws.Range("L132").Value = 6042
$ws.Range("M132").Value = -3884.529500000001
$ws.Range("N132").Value = -11102

$ws.Range("H134").Value = 2978.2666
$ws.Range("I134").Value = 2476.7144
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 7430.1432
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -4895.1432
$ws.Range("N134").Value = -35070

$ws.Range("H136").Value = 1324
$ws.Range("I136").Value = 1286.5
$ws.Range("J136").Value = 1499
$ws.Range("K136").Value = 3859.5
$ws.Range("L136").Value = 4497
$ws.Range("M136").Value = -1309.5
$ws.Range("N136").Value = -9597

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3157.5557
$ws.Range("J137").Value = 6165
$ws.Range("L137").Value = 18495
$ws.Range("N137").Value = -28695

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4515

$ws.Range("H46").Value = 45000
$ws.Range("J46").Value = 45000
$ws.Range("L46").Value = 45000
$ws.Range("N46").Value = -45312

$ws.Range("H80").Value = 6209.625
$ws.Range("I80").Value = 3890.3333
$ws.Range("J80").Value = 7601.2
$ws.Range("K80").Value = 3890.3333
$ws.Range("L80").Value = 7601.2
$ws.Range("M80").Value = -2892.3333
$ws.Range("N80").Value = -9597.200000000001

$ws.Range("H83").Value = 6209.625
$ws.Range("I83").Value = 3890.3333
$ws.Range("J83").Value = 7601.2
$ws.Range("K83").Value = 19451.6665
$ws.Range("L83").Value = 38006
$ws.Range("M83").Value = -14459.6665
$ws.Range("N83").Value = -47990

$ws.Range("H107").Value = 1423.35
$ws.Range("I107").Value = 256.33334
$ws.Range("J107").Value = 1923.5
$ws.Range("K107").Value = 256.33334
$ws.Range("L107").Value = 1923.5
$ws.Range("M107").Value = 1663.66666
$ws.Range("N107").Value = -5763.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6709.9
$ws.Range("J68").Value = 7478.7144
$ws.Range("L68").Value = 7478.7144
$ws.Range("N68").Value = -8976.714400000001

$ws.Range("H71").Value = 6709.9
$ws.Range("J71").Value = 7478.7144
$ws.Range("L71").Value = 37393.572
$ws.Range("N71").Value = -44881.572

$ws.Range("H82").Value = 2872.6875
$ws.Range("I82").Value = 1485.6666
$ws.Range("K82").Value = 1485.6666
$ws.Range("M82").Value = -1124.6666

$ws.Range("H85").Value = 2872.6875
$ws.Range("I85").Value = 1485.6666
$ws.Range("K85").Value = 1485.6666
$ws.Range("M85").Value = -237.6666

$ws.Range("H132").Value = 2958.3572
$ws.Range("I132").Value = 2898.5134
$ws.Range("J132").Value = 3401.2
$ws.Range("K132").Value = 8695.540199999999
$ws.Range("L132").Value = 10203.6
$ws.Range("M132").Value = -6165.540199999999
$ws.Range("N132").Value = -15263.6

$ws.Range("H136").Value = 2824.1353
$ws.Range("I136").Value = 1916.7778
$ws.Range("J136").Value = 3115.7856
$ws.Range("K136").Value = 5750.3334
$ws.Range("L136").Value = 9347.356800000001
$ws.Range("M136").Value = -3200.3334
$ws.Range("N136").Value = -14447.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15205.357
$ws.Range("J62").Value = 15472.223
$ws.Range("L62").Value = 15472.223
$ws.Range("N62").Value = -16720.223

$ws.Range("H65").Value = 15205.357
$ws.Range("J65").Value = 15472.223
$ws.Range("L65").Value = 77361.11500000001
$ws.Range("N65").Value = -83601.11500000001

$ws.Range("H81").Value = 3307.6875
$ws.Range("I81").Value = 2462
$ws.Range("J81").Value = 4153.375
$ws.Range("K81").Value = 4924
$ws.Range("L81").Value = 8306.75
$ws.Range("M81").Value = -3863
$ws.Range("N81").Value = -10428.75

$ws.Range("H84").Value = 3307.6875
$ws.Range("I84").Value = 2462
$ws.Range("J84").Value = 4153.375
$ws.Range("K84").Value = 24620
$ws.Range("L84").Value = 41533.75
$ws.Range("M84").Value = -19316
$ws.Range("N84").Value = -52141.75

$ws.Range("H113").Value = 1683.2222
$ws.Range("I113").Value = 1230.8462
$ws.Range("K113").Value = 3692.5386
$ws.Range("M113").Value = -1522.5386

$ws.Range("H132").Value = 10619.654
$ws.Range("I132").Value = 14655.823
$ws.Range("J132").Value = 2995.7778
$ws.Range("K132").Value = 43967.469
$ws.Range("L132").Value = 8987.3334
$ws.Range("M132").Value = -41437.469
$ws.Range("N132").Value = -14047.3334

$ws.Range("H136").Value = 1696.5834
$ws.Range("I136").Value = 1441.7273
$ws.Range("K136").Value = 4325.1819
$ws.Range("M136").Value = -1775.1819
